# Update timestamps on the "data" sheet (F2:F6)
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:33:13.961549"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:13.961557"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:13.961561"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:13.961564"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:13.961566"

# Add new "metadata" sheet after "data"
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Apply header style (style index 1 used for headers in "data" sheet)
$metaSheet.Range("B1:G1").Style = $dataSheet.Range("B1").Style

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("A2").Style = $dataSheet.Range("A2").Style
$metaSheet.Range("B2").Value = "Atrial Fibrillation"
$metaSheet.Range("C2").Value = 210
$metaSheet.Range("D2").Value = "0.7"
$metaSheet.Range("E2").Value = "2021-09-02T07:51:49.560503Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:13.957829"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/210/?format=json"
